$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "audit_usefulness" column (G) entirely; everything to its
# right shifts one column to the left.
$ws.Range("G:G").EntireColumn.Delete()

# After that shift, the old "w_audit_usefulness" column (which used to be
# M) now sits at column L; remove it too.
$ws.Range("L:L").EntireColumn.Delete()

# Update the scores / values that changed for the gpt-4o row.
$ws.Range("B2").Value = 3
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.75
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 0.4
$ws.Range("J2").Value = 0.3
$ws.Range("K2").Value = 0.9
$ws.Range("L2").Value = 3.35
$ws.Range("M2").Value = "The report provides a structured overview of debiasing strategies with clear sections and tables, earning a 4 in structure and formatting. However, the evidence extraction quality is rated 3 due to a lack of full sentence quotations and some missing validation details. Coverage of debiasing dimensions is strong, with a variety of methods discussed, but some common strategies like resampling are not explicitly documented, leading to a score of 4. Relevance and faithfulness are reasonable, but some claims lack direct support from the sources, resulting in a 3. Missing disclosures are identified, but not comprehensively, also scoring a 3. Overall, the report is useful for audits but could benefit from more detailed validation evidence and explicit documentation of all methods."
